# Commit: "updated powerpoint slide to slide 6"
#
# Semantic changes applied to the deck:
#   1. Slide 2 ("What is salmans?" title) - the three separate runs
#      ("What is " / "salmans" / "?") are collapsed into a single run
#      by re-assigning the TextRange.Text of the whole paragraph/shape.
#      Re-setting through a different intermediate value forces the
#      engine to rebuild the run list (rather than no-op when the
#      concatenated text is unchanged), exactly matching how PowerPoint
#      merges runs on a text assignment.
#   2. Slide 3 (picture shape) - a negligible (1 EMU) width correction
#      picked up the next time the picture's size was nudged/resaved by
#      PowerPoint.  We reproduce the exact resulting EMU value by setting
#      Width to the equivalent point value.

$p = $ppt.ActivePresentation

# --- Slide 2: merge "What is " + "salmans" + "?" runs into one run ---
$s2 = $p.Slides.Item(2)
$titleShape = $s2.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "temp"
$titleShape.TextFrame.TextRange.Text = "What is salmans?"

# --- Slide 3: picture width nudged from 8807115 EMU to 8807114 EMU ---
$s3 = $p.Slides.Item(3)
$picShape = $s3.Shapes.Item(2)
$picShape.Width = 693.4735433070866
